$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.5747214936203638
    3  = 0.5747214936203638
    4  = 0.4573291322971118
    5  = 0.4573291322971118
    6  = 0.4573291322971118
    7  = 0.9099502856655908
    8  = 0.6689697604315441
    9  = 0.8326105656603114
    10 = 0.8786003547184101
    11 = 0.8786003547184101
    12 = 0.8786003547184101
    13 = 1.386957062992251
    14 = 1.692375571059443
    15 = 1.546870014982985
    16 = 1.546870014982985
    17 = 4.872880388587285
    18 = 4.872880388587285
    19 = 4.872880388587285
    20 = 10.16437674197399
    21 = 10.16437674197399
}

foreach ($row in $values.Keys) {
    $ws.Range("Q$row").Value = $values[$row]
}

$wb.Save()
